$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.544.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.34%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.487.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.78%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.67'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.515.40'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.61%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.37'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.358'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.956.19'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.38'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.521.72'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.508.20'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.51'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.37'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.13%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.80'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.446'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -10.14%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.62%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.610.96'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.84'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.15'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0787'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.28'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.84'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.29%  '

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.39'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.42'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.73'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.50'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.66'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.98'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '315.36'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.69'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.76'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.837'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.80'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0940'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.21'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0528'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0231'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.22%  '
